$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.021.58"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "1.640.00"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.18"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5054"
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.008"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2575"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06446"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.47"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "1.648.64"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.244"
$ws.Range("D14").Value = "1.864.70"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5444"
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").Value = "0.0₅7910"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.48"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").Value = "26.001.25"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.99"
$ws.Range("E20").Value = "  -3.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.282"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.985"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.955"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.939"
$ws.Range("E25").Value = "  +10.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.78"
$ws.Range("E26").Value = "  -1.62%  "
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.72"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05050"
$ws.Range("E30").Value = "  -4.29%  "
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.247"
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.193"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.540"
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.636"
$ws.Range("E36").Value = "  -4.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8884"
$ws.Range("E37").Value = "  -3.96%  "
$ws.Range("D38").Value = "1.151.04"
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5612"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01572"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.564"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.666"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8071"
$ws.Range("E44").Value = "  -3.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.75"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").Value = "1.776.61"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4529"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.92"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05038"
$ws.Range("E51").Value = "  -0.85%  "
